$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = [double]"25.63000000000057"
$ws.Range("G2").Value = [double]"1.854072451124011e-14"
$ws.Range("H2").Value = [double]"2.849572429871235e-13"
$ws.Range("K2").Value = [double]"39.86990529515426"
$ws.Range("L2").Value = "[28.685230265483504, 51.05458032482502]"
$ws.Range("M2").Value = [double]"3.622258049063021e-11"
$ws.Range("N2").Value = [double]"3.622258049063021e-11"
$ws.Range("O2").Value = [double]"1.842816111114733"
$ws.Range("P2").Value = "[1.540921321580579, 2.144710900648888]"
$ws.Range("S2").Value = [double]"59.74057958255761"
$ws.Range("T2").Value = "[53.25549793978999, 66.22566122532523]"
$ws.Range("W2").Value = [double]"18.11289289289329"
$ws.Range("X2").Value = [double]"16.88142142142179"
$ws.Range("Y2").Value = [double]"19.34436436436479"

# Row 3 updates
$ws.Range("E3").Value = [double]"24.87000000000045"
$ws.Range("G3").Value = [double]"5.450084827884893e-13"
$ws.Range("H3").Value = [double]"2.515265482805224e-12"
$ws.Range("K3").Value = [double]"37.62971967224194"
$ws.Range("L3").Value = "[27.410584190997007, 47.84885515348687]"
$ws.Range("M3").Value = [double]"9.584999460798826e-12"
$ws.Range("N3").Value = [double]"1.916999892159765e-11"
$ws.Range("O3").Value = [double]"1.641552918091964"
$ws.Range("P3").Value = "[1.3396581285578097, 1.9434477076261185]"
$ws.Range("Q3").Value = [double]"0"
$ws.Range("R3").Value = [double]"0"
$ws.Range("S3").Value = [double]"57.21007340685843"
$ws.Range("T3").Value = "[50.832684747468015, 63.587462066248854]"
$ws.Range("W3").Value = [double]"18.37243243243276"
$ws.Range("X3").Value = [double]"17.17747747747779"
$ws.Range("Y3").Value = [double]"19.56738738738774"
